$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.399.67'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +3.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.836.91'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +3.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.026'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +2.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.46%  '
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('E7').Value = '  +2.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3718'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07343'
$ws.Range('D9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8715'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.32'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.904.31'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.88%  '
$ws.Range('E13').Value = '  +4.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.685'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07117'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.20'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.36%  '
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008979'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.47%  '
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E20').Value = '  +3.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.424.61'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.244'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.14'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '156.71'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.93%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.903'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.04%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.53'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.10%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.239'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.35%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.923'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +7.83%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '115.62'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09038'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('B31').Value = 'ARBITRUM'
$ws.Range('C31').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.198'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +7.54%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7586'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.88%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.463'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.58%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.862'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.025'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.10%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.143'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.65%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01956'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.91%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05247'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5155'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.03%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.785'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.03%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1661'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.99%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.540'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.10%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.465'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +6.49%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '108.35'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.50%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.50'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.025'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.16%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.674'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.20%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4618'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.43%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06299'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.870'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +9.20%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.38'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.16%  '
